# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.403.54"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").Value = "'1.694.39"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.26%  '

$ws.Range("D4").Value = "'1.010"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("D5").Value = "'219.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("D6").Value = "'0.5487"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.34%  '

$ws.Range("E7").Value = '  +0.30%  '

$ws.Range("D8").Value = "'0.2737"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.25%  '

$ws.Range("D9").Value = "'0.06468"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.39%  '

$ws.Range("E10").Value = '  -0.14%  '

$ws.Range("D11").Value = "'0.07669"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.64%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = "'4.555"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = "'1.674.04"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.92%  '

$ws.Range("D14").Value = "'0.5845"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.16%  '

$ws.Range("D15").Value = "'0.000008374"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").Value = "'65.44"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("D17").Value = "'26.456.66"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").Value = "'4.944"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.29%  '

$ws.Range("D19").Value = "'1.010"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.33%  '

$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("D21").Value = "'192.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.34%  '

$ws.Range("D22").Value = "'6.258"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.61%  '

$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").Value = "'149.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.96%  '

$ws.Range("D25").Value = "'0.1328"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.62%  '

$ws.Range("D26").Value = "'7.920"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.29%  '

$ws.Range("E27").Value = '  -0.62%  '

$ws.Range("D28").Value = "'0.06290"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -5.91%  '

$ws.Range("D29").Value = "'1.388"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.45%  '

$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("D32").Value = "'3.611"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.18%  '

$ws.Range("D33").Value = "'1.686"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.85%  '

$ws.Range("E34").Value = '  +1.50%  '

$ws.Range("D35").Value = "'0.6145"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.28%  '

$ws.Range("D36").Value = "'2.409"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("D37").Value = "'2.710"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.44%  '

$ws.Range("E38").Value = '  -2.49%  '

$ws.Range("D39").Value = "'0.01641"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.03%  '

$ws.Range("D40").Value = "'1.118.27"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.03%  '

$ws.Range("D41").Value = "'0.8888"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.50%  '

$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("D43").Value = "'101.91"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.08%  '

$ws.Range("D44").Value = "'1.845.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.39%  '

$ws.Range("D45").Value = "'0.00000000109"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.67%  '

$ws.Range("D46").Value = "'57.56"
$ws.Range("D46").ClearFormats()

$ws.Range("D47").Value = "'8.191"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("E48").Value = '  -0.25%  '

$ws.Range("D49").Value = "'0.05284"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.34%  '

$ws.Range("D50").Value = "'0.4304"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.06%  '

$ws.Range("D51").Value = "'6.096"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.61%  '
